# New crime data collected
# This script updates the 111th Precinct weekly CompStat report:
#  - bumps the "Volume ... Number" run from 2 -> 3
#  - shifts the reporting week dates forward by one week
#  - refreshes the crime-count/percentage figures for rows 16-21, 24-25, 27
#    (and a couple of cells that flip between numeric 0 counts and the
#    "***.*" / "0" placeholder text used elsewhere on the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Shared-string text tweaks (header runs)
# ---------------------------------------------------------------------
$ws.Range("C8").Characters(14, 1).Text = "3"

$ws.Range("C9").Characters(27, 7).Text = "1/16/2023"
$ws.Range("C9").Characters(47, 8).Text = "1/22/2023"

# ---------------------------------------------------------------------
# 2. Helper: reusable "placeholder" source cells. C14 already holds the
#    shared string "0" rendered with style 14, and E14 already holds the
#    shared string "***.*" rendered with style 14. D16 is a plain numeric
#    cell using style 16. Copying one of these onto a destination (via
#    Range.Copy(destination)) brings across both the value AND the exact
#    cell style without creating any unused style entries; we then
#    overwrite the value with the real figure we need.
# ---------------------------------------------------------------------
$zeroText = $ws.Range("C14")     # style 14, shared string "0"
$naText   = $ws.Range("E14")     # style 14, shared string "***.*"
$numStyle = $ws.Range("D16")     # style 16, plain number

# ---------------------------------------------------------------------
# 3. Row 15 : F15 numeric 1 -> text "0"
# ---------------------------------------------------------------------
$zeroText.Copy($ws.Range("F15"))
$ws.Range("F15").Value = "0"

# ---------------------------------------------------------------------
# 4. Row 16
# ---------------------------------------------------------------------
$zeroText.Copy($ws.Range("C16"))
$ws.Range("C16").Value = "0"

$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -33.333333333333
$ws.Range("J16").Value = 6
$ws.Range("K16").Value = -33.333333333333
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -76.470588235294

# ---------------------------------------------------------------------
# 5. Row 17 : C17 text "0" -> numeric 1
# ---------------------------------------------------------------------
$numStyle.Copy($ws.Range("C17"))
$ws.Range("C17").Value = 1

$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 4
$ws.Range("L17").Value = -33.333333333333
$ws.Range("N17").Value = -77.777777777777

# ---------------------------------------------------------------------
# 6. Row 18
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 30
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = 47.058823529411
$ws.Range("L18").Value = 19.047619047619
$ws.Range("M18").Value = 47.058823529411
$ws.Range("N18").Value = -60.9375

# ---------------------------------------------------------------------
# 7. Row 19
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -21.428571428571
$ws.Range("F19").Value = 57
$ws.Range("H19").Value = 5.555555555555
$ws.Range("I19").Value = 41
$ws.Range("J19").Value = 38
$ws.Range("K19").Value = 7.894736842105
$ws.Range("L19").Value = 173.333333333333
$ws.Range("M19").Value = 64
$ws.Range("N19").Value = 46.428571428571

# ---------------------------------------------------------------------
# 8. Row 20
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 25
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = 57.142857142857
$ws.Range("L20").Value = 266.666666666667
$ws.Range("M20").Value = 120
$ws.Range("N20").Value = -94.977168949771

# ---------------------------------------------------------------------
# 9. Row 21
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 3.333333333333
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = 10.309278350515
$ws.Range("I21").Value = 83
$ws.Range("J21").Value = 72
$ws.Range("K21").Value = 15.277777777777
$ws.Range("L21").Value = 97.619047619047
$ws.Range("M21").Value = 45.614035087719
$ws.Range("N21").Value = -75.443786982248

# ---------------------------------------------------------------------
# 10. Row 24
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 7
$ws.Range("E24").Value = -53.333333333333
$ws.Range("F24").Value = 39
$ws.Range("H24").Value = -41.791044776119
$ws.Range("I24").Value = 27
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = -46
$ws.Range("L24").Value = -20.588235294117
$ws.Range("M24").Value = -15.625

# ---------------------------------------------------------------------
# 11. Row 25
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 18
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 17
$ws.Range("J25").Value = 14
$ws.Range("K25").Value = 21.428571428571
$ws.Range("L25").Value = 183.333333333333
$ws.Range("M25").Value = 41.666666666666

# ---------------------------------------------------------------------
# 12. Row 26 : F26 numeric 1 -> text "0"
# ---------------------------------------------------------------------
$zeroText.Copy($ws.Range("F26"))
$ws.Range("F26").Value = "0"

# ---------------------------------------------------------------------
# 13. Row 27
# ---------------------------------------------------------------------
$ws.Range("F27").Value = 2

$zeroText.Copy($ws.Range("G27"))
$ws.Range("G27").Value = "0"

$naText.Copy($ws.Range("H27"))
$ws.Range("H27").Value = "***.*"

$ws.Range("I27").Value = 2
$ws.Range("L27").Value = 100
